$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 624, pushing existing rows 624:667 down to 625:668.
$ws.Rows.Item(624).Insert()

# Populate the newly inserted row 624 with the new record.
$ws.Cells.Item(624, 1).Value = 10
$ws.Cells.Item(624, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(624, 3).Value = "La Araucanía"
$ws.Cells.Item(624, 4).Value = 44931
$ws.Cells.Item(624, 5).Value = 9
$ws.Cells.Item(624, 6).Value = 100112032
$ws.Cells.Item(624, 7).Value = "Zapallo italiano"
$ws.Cells.Item(624, 8).Value = "Sin especificar"
$ws.Cells.Item(624, 9).Value = "Primera"
$ws.Cells.Item(624, 10).Value = 350
$ws.Cells.Item(624, 11).Value = 10000
$ws.Cells.Item(624, 12).Value = 10000
$ws.Cells.Item(624, 13).Value = 10000
$ws.Cells.Item(624, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(624, 15).Value = "Región del Maule"
$ws.Cells.Item(624, 16).Value = 200
$ws.Cells.Item(624, 17).Value = 50
$ws.Cells.Item(624, 18).Value = "Hortaliza"
